$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "D" column/row header and data
$ws.Range("D1").Value = "D"
$ws.Range("D2").Value = 0.2
$ws.Range("D3").Value = 0.2

$ws.Range("A4").Value = "D"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 1

# Remove the special number-format style that was previously on C2
$ws.Range("C2").Style = "Normal"

# Update the selection to match the target state
$ws.Range("C6").Select() | Out-Null
